$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values, re-pulled from source data.
$updates = @{
    2  = 4
    3  = 3
    4  = 2
    5  = -3
    6  = 4
    7  = -4
    8  = -3
    9  = -4
    10 = -2
    11 = -2
    12 = 3
    13 = -2
    14 = -3
    15 = 2
    17 = 5
    18 = -6
    19 = -5
    20 = -9
    21 = -4
    22 = -1
    24 = -1
    26 = 4
    27 = 3
    28 = 2
    30 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
